$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3572
$ws.Range("F5").Value = 3572
$ws.Range("F6").Value = 260
$ws.Range("F7").Value = 5097
$ws.Range("F8").Value = 5097
$ws.Range("F9").Value = 523
$ws.Range("F10").Value = 356
$ws.Range("F11").Value = 200
$ws.Range("F16").Value = 697
$ws.Range("F17").Value = 315
$ws.Range("F23").Value = 4908
$ws.Range("F24").Value = 4908
$ws.Range("F28").Value = 6036
$ws.Range("F32").Value = 342
$ws.Range("F33").Value = 711
$ws.Range("F34").Value = 4444
$ws.Range("F36").Value = 122
$ws.Range("F38").Value = 1016
$ws.Range("F42").Value = 871
$ws.Range("F43").Value = 987

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1118

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1118
$ws.Range("F7").Value = 3572
$ws.Range("F8").Value = 3572
$ws.Range("F9").Value = 260
$ws.Range("F10").Value = 5097
$ws.Range("F11").Value = 5097
$ws.Range("F12").Value = 523
$ws.Range("F13").Value = 356
$ws.Range("F14").Value = 200
$ws.Range("F19").Value = 697
$ws.Range("F20").Value = 315
$ws.Range("F27").Value = 4908
$ws.Range("F28").Value = 4908
$ws.Range("F32").Value = 6036
$ws.Range("F36").Value = 342
$ws.Range("F37").Value = 711
$ws.Range("F38").Value = 4444
$ws.Range("F41").Value = 122
$ws.Range("F43").Value = 1016
$ws.Range("F47").Value = 871
$ws.Range("F48").Value = 987
